$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.6749969208881521
$ws1.Range("C2").Value = 0.6173879758738527
$ws1.Range("B3").Value = 0.4590978704780589
$ws1.Range("C3").Value = -1.319022437466081
$ws1.Range("B4").Value = 1.087281290452234
$ws1.Range("C4").Value = -0.9204998865577018

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.8537523546243585
$ws2.Range("C2").Value = -0.1827437608322265
$ws2.Range("B3").Value = 1.051795902961459
$ws2.Range("C3").Value = 0.1052155019571581
$ws2.Range("B4").Value = -2.011587255793011
$ws2.Range("C4").Value = 0.3131885082383439
